# Build v2.1.2: Fix SearchCriteria variants and Schemas sheet grouping/sorting
#
# For the "Body" sheet and the "200"/"400" response sheets, the old
# per-field documentation rows (dateTime, settlementBIC, errorCode, ...)
# are collapsed into a single "schema" row that just references the
# named schema (currentPosition.211207Request / ...Response / errorResponse).
# For the response sheets that only had the header rows (204, 401, 403,
# 404, 429, 500) a new row 3 is appended that references the relevant
# schema (currentPosition.211207Response for 204, errorResponse1 for the
# pure error codes).

$wb = $excel.ActiveWorkbook

function Set-SchemaRow {
    param($ws, $section, $name, $schemaName)

    $ws.Range("A3").Value = $section
    $ws.Range("B3").Value = $name
    $ws.Range("C3").Value = ""
    $ws.Range("D3").Value = ""
    $ws.Range("E3").Value = "schema"
    $ws.Range("F3").Value = ""
    $ws.Range("G3").Value = $schemaName
    $ws.Range("H3").Value = ""
    $ws.Range("I3").Value = "Yes"
    $ws.Range("J3").Value = ""
    $ws.Range("K3").Value = ""
    $ws.Range("L3").Value = ""
    $ws.Range("M3").Value = ""
    $ws.Range("N3").Value = ""
    $ws.Range("O3").Value = ""
}

# --- Body sheet: request body now just references the request schema ---
$wsBody = $wb.Worksheets.Item("Body")
Set-SchemaRow $wsBody "body" "currentPosition.211207Request" "currentPosition.211207Request"
$wsBody.Rows("4").Delete()

# --- 200 response: collapse to a reference to the response schema ---
$ws200 = $wb.Worksheets.Item("200")
Set-SchemaRow $ws200 "content" "currentPosition.211207Response" "currentPosition.211207Response"
$ws200.Rows("4:7").Delete()

# --- 204 response: add a row referencing the response schema ---
$ws204 = $wb.Worksheets.Item("204")
Set-SchemaRow $ws204 "content" "currentPosition.211207Response" "currentPosition.211207Response"

# --- 400 response: collapse to a reference to the base error schema ---
$ws400 = $wb.Worksheets.Item("400")
Set-SchemaRow $ws400 "content" "errorResponse" "errorResponse"
$ws400.Rows("4:6").Delete()

# --- 401/403/404/429/500 responses: add a row referencing errorResponse1 ---
$errorCodeSheets = @("401", "403", "404", "429", "500")
foreach ($code in $errorCodeSheets) {
    $ws = $wb.Worksheets.Item($code)
    Set-SchemaRow $ws "content" "errorResponse1" "errorResponse1"
}
